$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Planilha de Testes"
